# Automatische test-sync: 2025-08-18 21:41:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append new row 13 to the Logs sheet
$logs.Range("A13").Value = "Geen onderwerp"
$logs.Range("B13").Value = "no-reply@testbedrijf123.nl"
$logs.Range("D13").Value = "Overig"
$logs.Range("F13").Value = "2025-08-18 21:40:55"
$logs.Range("G13").Value = "Nee"
$logs.Range("H13").Value = "Ja"
$logs.Range("I13").Value = "Nee"
$logs.Range("J13").Value = "Nee"

# Extend the conditional formatting ranges to include the new row (13)
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range($col + "2:" + $col + "12")
    $newRange = $logs.Range($col + "2:" + $col + "13")
    $count = $oldRange.FormatConditions.Count
    for ($i = 1; $i -le $count; $i++) {
        $fc = $oldRange.FormatConditions.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Overig"
$dashboard.Range("B3").Value = 5
